# Update "want to go" counts (column F) on the 展览, 演出 and 全部类型 sheets
# to reflect newly scraped numbers (gh-pages data refresh commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1549
$ws1.Range("F3").Value  = 8858
$ws1.Range("F4").Value  = 98
$ws1.Range("F5").Value  = 498
$ws1.Range("F6").Value  = 663
$ws1.Range("F7").Value  = 319
$ws1.Range("F10").Value = 42
$ws1.Range("F11").Value = 3741
$ws1.Range("F12").Value = 55
$ws1.Range("F13").Value = 373
$ws1.Range("F14").Value = 96
$ws1.Range("F15").Value = 3698
$ws1.Range("F18").Value = 323
$ws1.Range("F19").Value = 228
$ws1.Range("F20").Value = 2527
$ws1.Range("F21").Value = 89

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 37

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1549
$ws4.Range("F3").Value  = 8858
$ws4.Range("F4").Value  = 98
$ws4.Range("F5").Value  = 498
$ws4.Range("F6").Value  = 663
$ws4.Range("F7").Value  = 319
$ws4.Range("F10").Value = 42
$ws4.Range("F11").Value = 3741
$ws4.Range("F12").Value = 55
$ws4.Range("F13").Value = 373
$ws4.Range("F14").Value = 96
$ws4.Range("F15").Value = 3698
$ws4.Range("F18").Value = 323
$ws4.Range("F19").Value = 228
$ws4.Range("F20").Value = 2527
$ws4.Range("F21").Value = 37
$ws4.Range("F22").Value = 89

$wb.Save()
